$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B9").Value = "all length are forklength except for larvae is standard length"
$ws.Range("B10").Select() | Out-Null
